$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "country" column before the existing "city" column (F) ---
# This shifts city/address/province/zip/phone from F:J to G:K and leaves F
# free for the new "country" field.
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "country"

# --- Complete row 2 (existing account) with the newly-tracked fields ---
$ws.Range("F2").Value = "Indonesia"
$ws.Range("G2").Value = "luxemburg"
$ws.Range("H2").Value = "this street"
$ws.Range("I2").Value = "Lampung"
$ws.Range("J2").Value = 35111
$ws.Range("K2").Value = "'+6281000008"

# Re-point the existing A2 hyperlink/email at the refreshed test account.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Value = "john9935@gmail.com"
$ws.Range("B2").Value = "john123"
$ws.Range("C2").Value = "john9935"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:john9935@gmail.com") | Out-Null
$ws.Range("A2").NumberFormat = "@"

# --- New dummy accounts (rows 3-8) used by the delete-function tests ---
$accounts = @(
    @{ Row = 3; Name = "john9936"; Zip = 35112; Phone = "+6281000009" },
    @{ Row = 4; Name = "john9937"; Zip = 35113; Phone = "+6281000008" },
    @{ Row = 5; Name = "john9938"; Zip = 35114; Phone = "+6281000009" },
    @{ Row = 6; Name = "john9939"; Zip = 35115; Phone = "+6281000008" },
    @{ Row = 7; Name = "john9940"; Zip = 35116; Phone = "+6281000009" },
    @{ Row = 8; Name = "john9941"; Zip = 35117; Phone = "+6281000008" }
)

foreach ($acct in $accounts) {
    $r = $acct.Row
    $email = "$($acct.Name)@gmail.com"

    $ws.Range("A$r").Value = $email
    $ws.Range("B$r").Value = "john123"
    $ws.Range("C$r").Value = $acct.Name
    $ws.Range("D$r").Value = "john"
    $ws.Range("E$r").Value = "doe"
    $ws.Range("F$r").Value = "Indonesia"
    $ws.Range("G$r").Value = "luxemburg"
    $ws.Range("H$r").Value = "this street"
    $ws.Range("I$r").Value = "Lampung"
    $ws.Range("J$r").Value = $acct.Zip
    $ws.Range("K$r").Value = "'" + $acct.Phone

    $ws.Hyperlinks.Add($ws.Range("A$r"), "mailto:$email") | Out-Null
    $ws.Range("A$r").NumberFormat = "@"
}

# --- Column formatting to match the refreshed layout ---
$ws.Range("A1").NumberFormat = "@"
$ws.Columns("F").ColumnWidth = 11.17
$ws.Columns("H").ColumnWidth = 12.6
$ws.Columns("K").ColumnWidth = 13.6

# --- Misc view / print tweaks ---
$ws.Range("P22").Select() | Out-Null
$ws.PageSetup.Orientation = 1
